# Add a new order row (row 11) for the "BEC" sensor, mirroring the
# existing rows: a name in column A, a hyperlinked "link" label in
# column C, and a price in column D using a new Euro currency format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A11: sensor name -------------------------------------------------
$ws.Range("A11").Value = "BEC"

# --- D11: price, using a new numeric format (mirrors the new numFmt
#          169 "[$€-2] #,##0.00;[Red]-[$€-2] #,##0.00") --------------
# Set this before touching the hyperlink cell so its new style slot is
# allocated cleanly (index 3), matching the authored style table.
$ws.Range("D11").Value = 2.68
$ws.Range("D11").NumberFormat = "[$€-2]\ #,##0.00;[Red]\-[$€-2]\ #,##0.00"

# --- C11: hyperlinked "link" label, like C4/C6/C9/C10 -----------------
$ws.Range("C11").Value = "link"
$ws.Hyperlinks.Add($ws.Range("C11"), "https://www.ebay.com/itm/BEC-sensor") | Out-Null

# Re-apply the same look as the other "link" cells (reuses the existing
# Hyperlink-derived cell style instead of leaving the ad-hoc style that
# Hyperlinks.Add stamps on by default).
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Match the saved selection/active cell.
$ws.Range("D11").Select() | Out-Null
